$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before D ("Estimate"), shifting Estimate/Acceptance/Sprint right.
$ws.Columns("D").Insert()

# New column header + subtasks detail text.
$ws.Range("D1").Value = "Subtasks"
$ws.Range("D5").Value = "* Set up IDE `n* Set up PHP framework`n* Create database table for called Users`n* Create registration form`n* Create query to add user information to the database table`n* Form input validation `n"

# Rows whose "User Story Details" cell uses the plain (non-bordered) style never
# had an Estimate-column neighbor with a distinct look; after the insert those
# rows should have no Subtasks cell at all (fully cleared, not just blanked).
$clearRows = @(8,11,14,15,16,18,19,20,21,22,24,25,26,28)
foreach ($r in $clearRows) {
    $ws.Range("D$r").Clear()
}

# Support-registration story's estimate grows from 3 to 7 once its subtasks are broken out.
$ws.Range("E5").Value = 7
